{"js": "// Update the division expressions in the practice table to the new\n// set of dividends/divisors, matching each cell by its exact old text.\nconst body = context.document.body;\n\nconst replacements = [\n  [\"176\u00f79=\", \"476\u00f79=\"],\n  [\"745\u00f75=\", \"471\u00f77=\"],\n  [\"387\u00f73=\", \"590\u00f77=\"],\n  [\"927\u00f72=\", \"986\u00f78=\"],\n  [\"546\u00f73=\", \"651\u00f74=\"],\n  [\"921\u00f72=\", \"510\u00f77=\"],\n  [\"572\u00f78=\", \"353\u00f77=\"],\n  [\"452\u00f74=\", \"251\u00f77=\"],\n  [\"962\u00f78=\", \"958\u00f73=\"],\n  [\"372\u00f72=\", \"147\u00f78=\"],\n  [\"426\u00f74=\", \"355\u00f73=\"],\n  [\"572\u00f73=\", \"869\u00f73=\"],\n  [\"403\u00f72=\", \"232\u00f76=\"],\n  [\"473\u00f75=\", \"940\u00f75=\"],\n  [\"145\u00f75=\", \"812\u00f73=\"],\n  [\"569\u00f76=\", \"844\u00f76=\"],\n  [\"418\u00f74=\", \"598\u00f75=\"],\n  [\"743\u00f76=\", \"634\u00f79=\"],\n  [\"385\u00f73=\", \"900\u00f74=\"],\n  [\"511\u00f75=\", \"542\u00f77=\"],\n  [\"687\u00f78=\", \"529\u00f79=\"],\n  [\"746\u00f77=\", \"956\u00f73=\"],\n  [\"277\u00f78=\", \"849\u00f74=\"],\n  [\"953\u00f72=\", \"968\u00f74=\"],\n  [\"414\u00f76=\", \"894\u00f72=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the three-digit-division-by-one-digit expressions in the table\n# with the updated set of numbers, matching each cell by its exact old text.\n$d = $word.ActiveDocument\n\n$replacements = [ordered]@{\n  \"176\u00f79=\" = \"476\u00f79=\"\n  \"745\u00f75=\" = \"471\u00f77=\"\n  \"387\u00f73=\" = \"590\u00f77=\"\n  \"927\u00f72=\" = \"986\u00f78=\"\n  \"546\u00f73=\" = \"651\u00f74=\"\n  \"921\u00f72=\" = \"510\u00f77=\"\n  \"572\u00f78=\" = \"353\u00f77=\"\n  \"452\u00f74=\" = \"251\u00f77=\"\n  \"962\u00f78=\" = \"958\u00f73=\"\n  \"372\u00f72=\" = \"147\u00f78=\"\n  \"426\u00f74=\" = \"355\u00f73=\"\n  \"572\u00f73=\" = \"869\u00f73=\"\n  \"403\u00f72=\" = \"232\u00f76=\"\n  \"473\u00f75=\" = \"940\u00f75=\"\n  \"145\u00f75=\" = \"812\u00f73=\"\n  \"569\u00f76=\" = \"844\u00f76=\"\n  \"418\u00f74=\" = \"598\u00f75=\"\n  \"743\u00f76=\" = \"634\u00f79=\"\n  \"385\u00f73=\" = \"900\u00f74=\"\n  \"511\u00f75=\" = \"542\u00f77=\"\n  \"687\u00f78=\" = \"529\u00f79=\"\n  \"746\u00f77=\" = \"956\u00f73=\"\n  \"277\u00f78=\" = \"849\u00f74=\"\n  \"953\u00f72=\" = \"968\u00f74=\"\n  \"414\u00f76=\" = \"894\u00f72=\"\n}\n\nforeach ($old in $replacements.Keys) {\n  $new = $replacements[$old]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Text = $old\n  $find.Replacement.ClearFormatting()\n  $find.Replacement.Text = $new\n  $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
